$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of rows 5-19 (file/name pairs), sorted by the new UUID-prefixed
# filename. This reflects the table becoming sortable/scrollable while the
# underlying identifiers were regenerated.
$data = @(
    @("0d82fd06-f085-43a4-95f6-c01c3df9712d_fila_8.png", "David Orlando Mena Valverd"),
    @("1018ddf8-cd49-406d-9b4f-0328785670ef_fila_2.png", "Isabella Dompe Estrada"),
    @("5cb00160-fe91-436d-95da-2ac81e3d03d3_fila_1.png", "Hotep Antonio Ruiz Lezama"),
    @("6ce36d1e-8d74-4797-9ba3-fa7ba90432a2_fila_4.png", "Yadder Fernando Torres"),
    @("7428429f-b1e6-4600-a2a3-ef420d7a093d_fila_11.png", "Eduardo Domingo Zeledon Merca"),
    @("7e5534c5-5bad-4946-968f-43a009d6e689_fila_15.png", "Orlando Mauricio Guevara"),
    @("873975f9-9332-47bd-9209-318e28da4d12_fila_12.png", "José Danilo Suárez"),
    @("8a6f138b-ebda-4d30-86d7-fbaba74d0a4d_fila_3.png", "Bryan Alexander Cano"),
    @("98dc1945-9968-4cc9-813b-2b342b3b07e3_fila_14.png", "Angel Isaac Alvarez Quiñonez"),
    @("a723f7c1-2614-48c9-bf41-524713caf936_fila_6.png", "Ronier Jose Rivera"),
    @("bd356b71-dd2e-4871-867a-6249232a5812_fila_10.png", "Abraham Silva Ampre"),
    @("caea43b6-8cb6-4ed7-a00c-f4ffd6f19172_fila_9.png", "Roman Alfonso Grios Boza"),
    @("e0927981-dbc4-4282-98a7-d738d0f0c4c0_fila_5.png", "Erick Espinoza"),
    @("e6f13450-e309-449c-98f4-e20561f9ac2a_fila_13.png", "Marlon Josue Gonzales Cano"),
    @("e8216967-2979-4117-b634-25ba5bfd9a25_fila_7.png", "Cristina Jozabed Carvajal")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
